$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header columns to snake_case machine-readable names
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# 2) Title-case the Spanish connector words (de/del/la/las/los/el/y) in state/
#    municipality names throughout the data rows
$titleCaseFixes = @(
    @('B8', 'Pabellón De Arteaga'),
    @('B9', 'Rincón De Romos'),
    @('B10', 'San Francisco De Los Romo'),
    @('B11', 'San José De Gracia'),
    @('B40', 'Amatenango De La Frontera'),
    @('B44', 'Benemérito De Las Américas'),
    @('B57', 'Comitán De Domínguez'),
    @('B82', 'Marqués De Comillas'),
    @('B83', 'Mazapa De Madero'),
    @('B89', 'Ocozocoautla De Espinosa'),
    @('B99', 'Salto De Agua'),
    @('B100', 'San Cristóbal De Las Casas'),
    @('B144', 'Coyame Del Sotol'),
    @('B155', 'Guadalupe Y Calvo'),
    @('B158', 'Hidalgo Del Parral'),
    @('B184', 'San Francisco De Borja'),
    @('B185', 'San Francisco De Conchos'),
    @('B186', 'San Francisco Del Oro'),
    @('B194', 'Valle De Zaragoza'),
    @('B227', 'San Juan De Sabinas'),
    @('B244', 'Villa De Álvarez'),
    @('A246', 'Ciudad De México'),
    @('B250', 'Cuajimalpa De Morelos'),
    @('B265', 'Coneto De Comonfort'),
    @('B279', 'Nombre De Dios'),
    @('B283', 'Pánuco De Coronado'),
    @('B290', 'San Juan De Guadalupe'),
    @('B291', 'San Juan Del Río'),
    @('B292', 'San Luis Del Cordero'),
    @('B293', 'San Pedro Del Gallo'),
    @('A303', 'Estado De México'),
    @('B303', 'Acambay De Ruíz Castañeda'),
    @('B306', 'Almoloya De Alquisiras'),
    @('B307', 'Almoloya De Juárez'),
    @('B308', 'Almoloya Del Río'),
    @('B315', 'Atizapán De Zaragoza'),
    @('B322', 'Chapa De Mota'),
    @('B328', 'Coacalco De Berriozábal'),
    @('B335', 'Ecatepec De Morelos'),
    @('B343', 'Ixtapan De La Sal'),
    @('B344', 'Ixtapan Del Oro'),
    @('B359', 'Naucalpan De Juárez'),
    @('B373', 'San Antonio La Isla'),
    @('B374', 'San Felipe Del Progreso'),
    @('B375', 'San Martín De Las Pirámides'),
    @('B377', 'San Simón De Guerrero'),
    @('B379', 'Soyaniquilpan De Juárez'),
    @('B389', 'Tenango Del Aire'),
    @('B390', 'Tenango Del Valle'),
    @('B404', 'Tlalnepantla De Baz'),
    @('B410', 'Valle De Bravo'),
    @('B411', 'Valle De Chalco Solidaridad'),
    @('B412', 'Villa De Allende'),
    @('B413', 'Villa Del Carbón'),
    @('B426', 'San Miguel De Allende'),
    @('B427', 'Apaseo El Alto'),
    @('B428', 'Apaseo El Grande'),
    @('B436', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B440', 'Jaral Del Progreso'),
    @('B448', 'Purísima Del Rincón'),
    @('B452', 'San Diego De La Unión'),
    @('B454', 'San Francisco Del Rincón'),
    @('B456', 'San Luis De La Paz'),
    @('B458', 'Santa Cruz De Juventino Rosas'),
    @('B460', 'Silao De La Victoria'),
    @('B465', 'Valle De Santiago'),
    @('B471', 'Acapulco De Juárez'),
    @('B474', 'Ajuchitlán Del Progreso'),
    @('B475', 'Alcozauca De Guerrero'),
    @('B479', 'Atenango Del Río'),
    @('B480', 'Atlamajalcingo Del Monte'),
    @('B482', 'Atoyac De Álvarez'),
    @('B483', 'Ayutla De Los Libres'),
    @('B486', 'Buenavista De Cuéllar'),
    @('B487', 'Chilapa De Álvarez'),
    @('B488', 'Chilpancingo De Los Bravo'),
    @('B489', 'Coahuayutla De José María Izazaga'),
    @('B494', 'Coyuca De Benítez'),
    @('B495', 'Coyuca De Catalán'),
    @('B499', 'Cuetzala Del Progreso'),
    @('B500', 'Cutzamala De Pinzón'),
    @('B506', 'Huitzuco De Los Figueroa'),
    @('B507', 'Iguala De La Independencia'),
    @('B509', 'Ixcateopan De Cuauhtémoc'),
    @('B510', 'Zihuatanejo De Azueta'),
    @('B512', 'La Unión De Isidoro Montes De Oca'),
    @('B515', 'Mártir De Cuilapan'),
    @('B528', 'Taxco De Alarcón'),
    @('B530', 'Técpan De Galeana'),
    @('B532', 'Tepecoacuilco De Trujano'),
    @('B534', 'Tixtla De Guerrero'),
    @('B538', 'Tlalixtaquilla De Maldonado'),
    @('B539', 'Tlapa De Comonfort'),
    @('B551', 'Agua Blanca De Iturbide'),
    @('B558', 'Atotonilco De Tula'),
    @('B559', 'Atotonilco El Grande'),
    @('B565', 'Cuautepec De Hinojosa'),
    @('B571', 'Huasca De Ocampo'),
    @('B575', 'Huejutla De Reyes'),
    @('B578', 'Jacala De Ledezma'),
    @('B585', 'Mineral De La Reforma'),
    @('B586', 'Mineral Del Chico'),
    @('B587', 'Mineral Del Monte'),
    @('B588', 'Mixquiahuala De Juárez'),
    @('B589', 'Molango De Escamilla'),
    @('B591', 'Nopala De Villagrán'),
    @('B592', 'Omitlán De Juárez'),
    @('B593', 'Pachuca De Soto'),
    @('B596', 'Progreso De Obregón'),
    @('B602', 'Santiago De Anaya'),
    @('B603', 'Santiago Tulantepec De Lugo Guerrero'),
    @('B607', 'Tenango De Doria'),
    @('B609', 'Tepehuacán De Guerrero'),
    @('B610', 'Tepeji Del Río De Ocampo'),
    @('B613', 'Tezontepec De Aldama'),
    @('B622', 'Tula De Allende'),
    @('B623', 'Tulancingo De Bravo'),
    @('B624', 'Villa De Tezontepec'),
    @('B628', 'Zacualtipán De Ángeles'),
    @('B629', 'Zapotlán De Juárez'),
    @('B634', 'Acatlán De Juárez'),
    @('B635', 'Ahualulco De Mercado'),
    @('B640', 'Atemajac De Brizuela'),
    @('B643', 'Atotonilco El Alto'),
    @('B645', 'Autlán De Navarro'),
    @('B651', 'Cañadas De Obregón'),
    @('B658', 'Concepción De Buenos Aires'),
    @('B659', 'Cuautitlán De García Barragán'),
    @('B668', 'Encarnación De Díaz'),
    @('B675', 'Huejuquilla El Alto'),
    @('B676', 'Ixtlahuacán De Los Membrillos'),
    @('B677', 'Ixtlahuacán Del Río'),
    @('B681', 'Jilotlán De Los Dolores'),
    @('B687', 'La Manzanilla De La Paz'),
    @('B688', 'Lagos De Moreno'),
    @('B695', 'Ojuelos De Jalisco'),
    @('B700', 'San Cristóbal De La Barranca'),
    @('B701', 'San Diego De Alejandría'),
    @('B703', 'San Juan De Los Lagos'),
    @('B706', 'San Martín De Bolaños'),
    @('B708', 'San Miguel El Alto'),
    @('B709', 'San Sebastián Del Oeste'),
    @('B710', 'Santa María De Los Ángeles'),
    @('B711', 'Santa María Del Oro'),
    @('B714', 'Talpa De Allende'),
    @('B715', 'Tamazula De Gordiano'),
    @('B718', 'Techaluta De Montenegro'),
    @('B722', 'Teocuitatlán De Corona'),
    @('B723', 'Tepatitlán De Morelos'),
    @('B726', 'Tizapán El Alto'),
    @('B727', 'Tlajomulco De Zúñiga'),
    @('B739', 'Unión De San Antonio'),
    @('B740', 'Unión De Tula'),
    @('B741', 'Valle De Guadalupe'),
    @('B742', 'Valle De Juárez'),
    @('B747', 'Yahualica De González Gallo'),
    @('B748', 'Zacoalco De Torres'),
    @('B751', 'Zapotitlán De Vadillo'),
    @('B752', 'Zapotlán Del Rey'),
    @('B753', 'Zapotlán El Grande'),
    @('B779', 'Coalcomán De Vázquez Pallares'),
    @('B781', 'Cojumatlán De Régules'),
    @('B848', 'Tiquicheo De Nicolás Romero'),
    @('B874', 'Coatlán Del Río'),
    @('B882', 'Jonacatepec De Leandro Valle'),
    @('B886', 'Puente De Ixtla'),
    @('B892', 'Tetela Del Volcán'),
    @('B893', 'Tlaltizapán De Zapata'),
    @('B901', 'Zacualpan De Amilpas'),
    @('B905', 'Amatlán De Cañas'),
    @('B906', 'Bahía De Banderas'),
    @('B910', 'Ixtlán Del Río'),
    @('B917', 'Santa María Del Oro'),
    @('B935', 'Ciénega De Flores'),
    @('B953', 'Lampazos De Naranjo'),
    @('B960', 'Mier Y Noriega'),
    @('B969', 'San Nicolás De Los Garza'),
    @('B977', 'Acatlán De Pérez Figueroa'),
    @('B983', 'Ayoquezco De Aldama'),
    @('B986', 'Capulálpam De Méndez'),
    @('B988', 'Chalcatongo De Hidalgo'),
    @('B989', 'Chiquihuitlán De Benito Juárez'),
    @('B992', 'Coicoyán De Las Flores'),
    @('B994', 'Constancia Del Rosario'),
    @('B996', 'El Barrio De La Soledad'),
    @('B998', 'Fresnillo De Trujano'),
    @('B999', 'Guadalupe De Ramírez'),
    @('B1001', 'Guevea De Humboldt'),
    @('B1002', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B1003', 'Heroica Ciudad De Huajuapan De León'),
    @('B1004', 'Heroica Ciudad De Tlaxiaco'),
    @('B1006', 'Ixtlán De Juárez'),
    @('B1007', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B1016', 'Magdalena Yodocono De Porfirio Díaz'),
    @('B1017', 'Mariscala De Juárez'),
    @('B1018', 'Mártires De Tacubaya'),
    @('B1021', 'Miahuatlán De Porfirio Díaz'),
    @('B1024', 'Nejapa De Madero'),
    @('B1026', 'Oaxaca De Juárez'),
    @('B1027', 'Ocotlán De Morelos'),
    @('B1028', 'Pinotepa De Don Luis'),
    @('B1030', 'Putla Villa De Guerrero'),
    @('B1031', 'Reforma De Pineda'),
    @('B1062', 'San Francisco Del Mar'),
    @('B1076', 'San José Del Peñasco'),
    @('B1077', 'San José Del Progreso'),
    @('B1091', 'San Juan De Los Cués'),
    @('B1092', 'San Juan Del Estado'),
    @('B1122', 'San Martín De Los Cansecos'),
    @('B1134', 'San Miguel Del Puerto'),
    @('B1136', 'San Miguel El Grande'),
    @('B1172', 'San Pedro Y San Pablo Teposcolula'),
    @('B1192', 'Santa Cruz Tacache De Mina'),
    @('B1196', 'Santa Inés De Zaragoza'),
    @('B1197', 'Santa Inés Del Monte'),
    @('B1204', 'Santa María Del Rosario'),
    @('B1205', 'Santa María Del Tule'),
    @('B1210', 'Santa María Jalapa Del Marqués'),
    @('B1211', 'Santa María La Asunción'),
    @('B1264', 'Santo Domingo De Morelos'),
    @('B1282', 'Tamazulápam Del Espíritu Santo'),
    @('B1284', 'Tataltepec De Valdés'),
    @('B1285', 'Teococuilco De Marcos Pérez'),
    @('B1286', 'Teotitlán De Flores Magón'),
    @('B1287', 'Teotitlán Del Valle'),
    @('B1288', 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'),
    @('B1289', 'Tlacolula De Matamoros'),
    @('B1290', 'Tlalixtac De Cabrera'),
    @('B1294', 'Villa De Chilapa De Díaz'),
    @('B1295', 'Villa De Etla'),
    @('B1296', 'Villa De Tamazulápam Del Progreso'),
    @('B1297', 'Villa De Tututepec'),
    @('B1298', 'Villa De Zaachila'),
    @('B1300', 'Villa Sola De Vega'),
    @('B1301', 'Villa Talea De Castro'),
    @('B1302', 'Zapotitlán Del Río'),
    @('B1304', 'Zimatlán De Álvarez'),
    @('B1327', 'Ayotoxco De Guerrero'),
    @('B1331', 'Chalchicomula De Sesma'),
    @('B1350', 'Cuapiaxtla De Madero'),
    @('B1354', 'Cuayuca De Andrade'),
    @('B1355', 'Cuetzalan Del Progreso'),
    @('B1369', 'Huehuetlán El Chico'),
    @('B1370', 'Huehuetlán El Grande'),
    @('B1373', 'Huitzilan De Serdán'),
    @('B1375', 'Ixcamilpa De Guerrero'),
    @('B1379', 'Izúcar De Matamoros'),
    @('B1389', 'Los Reyes De Juárez'),
    @('B1390', 'Mazapiltepec De Juárez'),
    @('B1400', 'Palmar De Bravo'),
    @('B1408', 'San Diego La Mesa Tochimiltzingo'),
    @('B1416', 'San Nicolás De Los Ranchos'),
    @('B1419', 'San Salvador El Seco'),
    @('B1420', 'San Salvador El Verde'),
    @('B1428', 'Tecali De Herrera'),
    @('B1436', 'Tepanco De López'),
    @('B1437', 'Tepatlaxco De Hidalgo'),
    @('B1443', 'Tepexi De Rodríguez'),
    @('B1445', 'Tepeyahualco De Cuauhtémoc'),
    @('B1446', 'Tetela De Ocampo'),
    @('B1447', 'Teteles De Avila Castillo'),
    @('B1452', 'Tlacotepec De Benito Juárez'),
    @('B1464', 'Totoltepec De Guerrero'),
    @('B1466', 'Tuzamapan De Galeana'),
    @('B1470', 'Xayacatlán De Bravo'),
    @('B1489', 'Amealco De Bonfil'),
    @('B1491', 'Cadereyta De Montes'),
    @('B1497', 'Jalpan De Serra'),
    @('B1498', 'Landa De Matamoros'),
    @('B1501', 'Pinal De Amoles'),
    @('B1504', 'San Juan Del Río'),
    @('B1517', 'Armadillo De Los Infante'),
    @('B1518', 'Axtla De Terrazas'),
    @('B1523', 'Cerro De San Pedro'),
    @('B1525', 'Ciudad Del Maíz'),
    @('B1536', 'Mexquitic De Carmona'),
    @('B1542', 'San Ciro De Acosta'),
    @('B1548', 'Santa María Del Río'),
    @('B1550', 'Soledad De Graciano Sánchez'),
    @('B1558', 'Tanquián De Escobedo'),
    @('B1562', 'Villa De Arista'),
    @('B1563', 'Villa De Arriaga'),
    @('B1564', 'Villa De Guadalupe'),
    @('B1565', 'Villa De La Paz'),
    @('B1566', 'Villa De Ramos'),
    @('B1567', 'Villa De Reyes'),
    @('B1613', 'Nacozari De García'),
    @('B1620', 'San Felipe De Jesús'),
    @('B1635', 'Jalpa De Méndez'),
    @('B1680', 'Soto La Marina'),
    @('B1688', 'Acuamanala De Miguel Hidalgo'),
    @('B1690', 'Amaxac De Guerrero'),
    @('B1691', 'Apetatitlán De Antonio Carvajal'),
    @('B1703', 'Ixtacuixtla De Mariano Matamoros'),
    @('B1706', 'Mazatecochco De José María Morelos'),
    @('B1707', 'Muñoz De Domingo Arenas'),
    @('B1708', 'Nanacamilpa De Mariano Arista'),
    @('B1711', 'Papalotla De Xicohténcatl'),
    @('B1714', 'San Pablo Del Monte'),
    @('B1715', 'Sanctórum De Lázaro Cárdenas'),
    @('B1721', 'Tetla De La Solidaridad'),
    @('B1732', 'Ziltlaltépec De Trinidad Sánchez Santos'),
    @('B1742', 'Alto Lucero De Gutiérrez Barrios'),
    @('B1746', 'Amatlán De Los Reyes'),
    @('B1758', 'Boca Del Río'),
    @('B1763', 'Castillo De Teayo'),
    @('B1765', 'Cazones De Herrera'),
    @('B1773', 'Chinampa De Gorostiza'),
    @('B1787', 'Cosamaloapan De Carpio'),
    @('B1788', 'Cosautlán De Carvajal'),
    @('B1805', 'Hueyapan De Ocampo'),
    @('B1806', 'Ignacio De La Llave'),
    @('B1810', 'Ixhuacán De Los Reyes'),
    @('B1811', 'Ixhuatlán De Madero'),
    @('B1812', 'Ixhuatlán Del Café'),
    @('B1813', 'Ixhuatlán Del Sureste'),
    @('B1824', 'Juchique De Ferrer'),
    @('B1827', 'Landero Y Coss'),
    @('B1830', 'Las Vigas De Ramírez'),
    @('B1831', 'Lerdo De Tejada'),
    @('B1836', 'Martínez De La Torre'),
    @('B1839', 'Medellín De Bravo'),
    @('B1843', 'Mixtla De Altamirano'),
    @('B1845', 'Nanchital De Lázaro Cárdenas Del Río'),
    @('B1856', 'Ozuluama De Mascareñas'),
    @('B1860', 'Paso De Ovejas'),
    @('B1861', 'Paso Del Macho'),
    @('B1865', 'Poza Rica De Hidalgo'),
    @('B1876', 'Sayula De Alemán'),
    @('B1879', 'Soledad De Doblado'),
    @('B1887', 'Tatahuicapan De Juárez'),
    @('B1908', 'Tlacotepec De Mejía'),
    @('B1922', 'Vega De Alatorre'),
    @('B1934', 'Zontecomatlán De López Y Fuentes'),
    @('B1935', 'Zozocolco De Hidalgo'),
    @('B1971', 'Cañitas De Felipe Pescador'),
    @('B1973', 'Concepción Del Oro'),
    @('B1975', 'El Plateado De Joaquín Amaro'),
    @('B1986', 'Jiménez Del Teul'),
    @('B1993', 'Mezquital Del Oro'),
    @('B1998', 'Moyahua De Estrada'),
    @('B1999', 'Nochistlán De Mejía'),
    @('B2000', 'Noria De Ángeles'),
    @('B2011', 'Teúl De González Ortega'),
    @('B2012', 'Tlaltenango De Sánchez Román'),
    @('B2014', 'Trinidad García De La Cadena'),
    @('B2017', 'Villa De Cos')
)
foreach ($fix in $titleCaseFixes) {
    $ws.Range($fix[0]).Value = $fix[1]
}

# 3) Apply the tiny (1-ULP) floating point recompute nudges on the percentage
#    column that Excel's own division produced for these particular totals
$numericFixes = @(
    @('D13', 0.009380353240002364),
    @('D158', 0.0009392167286904128),
    @('D303', 0.0009510307755921792),
    @('D426', 0.009498493709020023),
    @('D453', 0.009132258255065272),
    @('D456', 0.009474865615216491),
    @('D530', 0.0009628448224939452),
    @('D721', 0.0009510307755921792),
    @('D1646', 0.0009510307755921792),
    @('D1680', 0.0009628448224939452),
    @('D2017', 0.0009155886348868804)
)
foreach ($fix in $numericFixes) {
    $ws.Range($fix[0]).Value = $fix[1]
}

# 4) Remove the trailing footer/metadata rows (2026-2030) that are no longer part
#    of the clean data table
$ws.Range("A2026:A2030").EntireRow.Delete()
